# Update Excel workbook after daily scrape - 2025-08-25 03:24:04 UTC
# Row 2 is updated in place; rows 3-11 are newly appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments (C, D, G, H) -------------------------------
# ColumnWidth, as exposed by this COM surface, stores the value with a
# constant +5/6 character offset baked in by the engine's internal unit
# conversion. Subtract it up front so the persisted <col width="..."/>
# lands on the exact target values (30 / 65 / 16 / 60).
$colWidthOffset = 0.8333333333333334
$ws.Columns.Item(3).ColumnWidth = 30 - $colWidthOffset
$ws.Columns.Item(4).ColumnWidth = 65 - $colWidthOffset
$ws.Columns.Item(7).ColumnWidth = 16 - $colWidthOffset
$ws.Columns.Item(8).ColumnWidth = 60 - $colWidthOffset

# --- Mark column A (OPPORTUNITY ID) as Text before writing ---------------
# The IDs are purely numeric-looking strings ("1327071" etc.) and must stay
# text, matching the source feed / original workbook's inline strings.
# Pre-formatting the range as Text ("@") keeps Excel from coercing the
# written values to numbers.
$ws.Range("A2:A11").NumberFormat = "@"

# --- Row 2 (existing row, values replaced) + rows 3-11 (new rows) --------
$ws.Range("A2").Value = "1327071"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1327071"
$ws.Range("C2").Value = "Tax Intern"
$ws.Range("D2").Value = "Panamá, Provincia de Panamá, Panamá"
$ws.Range("E2").Value = "No"
$ws.Range("F2").Value = "2 applicants"
$ws.Range("G2").Value = "6 - 18 Months"
$ws.Range("H2").Value = "Samsung Electronics Latinoamérica (Zona Libre) S.A (SELA)"
$ws.Range("A3").Value = "1327067"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1327067"
$ws.Range("C3").Value = "Creative Marketing Intern"
$ws.Range("D3").Value = "Hyderabad, Telangana, India"
$ws.Range("E3").Value = "No"
$ws.Range("F3").Value = "0 applicants"
$ws.Range("G3").Value = "6 - 18 Months"
$ws.Range("H3").Value = "MPF clothing collection PVT LTD"
$ws.Range("A4").Value = "1326995"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1326995"
$ws.Range("C4").Value = "Sales Specialist"
$ws.Range("D4").Value = "10th of Ramadan City, Al-Sharqia Governorate, Egypt"
$ws.Range("E4").Value = "No"
$ws.Range("F4").Value = "1 applicant"
$ws.Range("G4").Value = "9 - 12 Weeks"
$ws.Range("H4").Value = "ABD Eldaem Road Transportation Company"
$ws.Range("A5").Value = "1326990"
$ws.Range("B5").Value = "https://aiesec.org/opportunity/global-talent/1326990"
$ws.Range("C5").Value = "Artificial intelligence"
$ws.Range("D5").Value = "10th of Ramadan City, Al-Sharqia Governorate, Egypt"
$ws.Range("E5").Value = "No"
$ws.Range("F5").Value = "4 applicants"
$ws.Range("G5").Value = "9 - 12 Weeks"
$ws.Range("H5").Value = "ABD Eldaem Road Transportation Company"
$ws.Range("A6").Value = "1326744"
$ws.Range("B6").Value = "https://aiesec.org/opportunity/global-talent/1326744"
$ws.Range("C6").Value = "Cyber Security - IT"
$ws.Range("D6").Value = "Maadi as Sarayat Al Gharbeyah, Maadi, Cairo Governorate, Egypt"
$ws.Range("E6").Value = "No"
$ws.Range("F6").Value = "2 applicants"
$ws.Range("G6").Value = "6 - 18 Months"
$ws.Range("H6").Value = "Keys Payroll"
$ws.Range("A7").Value = "1326743"
$ws.Range("B7").Value = "https://aiesec.org/opportunity/global-talent/1326743"
$ws.Range("C7").Value = "Marketing Specialist"
$ws.Range("D7").Value = "Maadi as Sarayat Al Gharbeyah, Maadi, Cairo Governorate, Egypt"
$ws.Range("E7").Value = "No"
$ws.Range("F7").Value = "1 applicant"
$ws.Range("G7").Value = "6 - 18 Months"
$ws.Range("H7").Value = "Keys Payroll"
$ws.Range("A8").Value = "1326741"
$ws.Range("B8").Value = "https://aiesec.org/opportunity/global-talent/1326741"
$ws.Range("C8").Value = "Business Development Intern"
$ws.Range("D8").Value = "Hyderabad, Telangana, India"
$ws.Range("E8").Value = "No"
$ws.Range("F8").Value = "0 applicants"
$ws.Range("G8").Value = "9 - 12 Weeks"
$ws.Range("H8").Value = "Vigilare biopharma Pvt Ltd"
$ws.Range("A9").Value = "1326713"
$ws.Range("B9").Value = "https://aiesec.org/opportunity/global-talent/1326713"
$ws.Range("C9").Value = "HR Specialist"
$ws.Range("D9").Value = "Maadi as Sarayat Al Gharbeyah, Maadi, Cairo Governorate, Egypt"
$ws.Range("E9").Value = "No"
$ws.Range("F9").Value = "0 applicants"
$ws.Range("G9").Value = "6 - 18 Months"
$ws.Range("H9").Value = "Keys Payroll"
$ws.Range("A10").Value = "1325656"
$ws.Range("B10").Value = "https://aiesec.org/opportunity/global-talent/1325656"
$ws.Range("C10").Value = "Design Intern"
$ws.Range("D10").Value = "Mumbai, Maharashtra, India"
$ws.Range("E10").Value = "No"
$ws.Range("F10").Value = "5 applicants"
$ws.Range("G10").Value = "6 - 18 Months"
$ws.Range("H10").Value = "Rediffusion Brand Solutions Pvt Ltd"
$ws.Range("A11").Value = "1325594"
$ws.Range("B11").Value = "https://aiesec.org/opportunity/global-talent/1325594"
$ws.Range("C11").Value = "Brand Partner"
$ws.Range("D11").Value = "Mumbai, Maharashtra, India"
$ws.Range("E11").Value = "No"
$ws.Range("F11").Value = "21 applicants"
$ws.Range("G11").Value = "6 - 18 Months"
$ws.Range("H11").Value = "Rediffusion Brand Solutions Pvt Ltd"

# --- Reset column A's style back to Normal --------------------------------
# Only the number FORMAT needed to be "@" at write-time to force text
# storage; clear the resulting style delta afterwards so the cells keep the
# workbook's default (unstyled) appearance, same as every other data cell.
$ws.Range("A2:A11").Style = "Normal"
